# Applies the "adding a template" edit to slide 1:
#  1. TextBox 9  - cosmetic "dirty" flag on the "40 horas" run (left as-is;
#                  the run's content/formatting is untouched).
#  2. Freeform 10 - right-align the (empty) paragraph in its text body.
#  3. AutoShape 14 - reposition/resize the decorative line, clearing the
#                    vertical flip.
#  4. AutoShape 15 - reposition/resize the decorative line, clearing the
#                    vertical flip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Freeform 10: right-align its paragraph -------------------------------
$logo = $s.Shapes.Item("Freeform 10")
$logo.TextFrame.TextRange.ParagraphFormat.Alignment = 3   # ppAlignRight

# --- AutoShape 14: move/resize the line, remove the vertical flip --------
$line1 = $s.Shapes.Item("AutoShape 14")
$line1.Flip(1)                      # msoFlipVertical - toggles flipV off
$line1.Left   = 125.506416321
$line1.Top    = 508.703979492
$line1.Width  = 227.780044556
$line1.Height = 0.659608245

# --- AutoShape 15: move/resize the line, remove the vertical flip --------
$line2 = $s.Shapes.Item("AutoShape 15")
$line2.Flip(1)                      # msoFlipVertical - toggles flipV off
$line2.Left   = 488.603363037
$line2.Top    = 509.363677979
$line2.Width  = 227.780044556
$line2.Height = 0.000000000
